$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update likes_count values for rows 5 and 6 (T column)
$ws.Range("T5").Value = 702
$ws.Range("T6").Value = 200

# Delete row 7 entirely (the_year_book_ / DJB1QwfPUKC post), shrinking the table
$ws.Rows.Item(7).Delete()

# Resize the table / autofilter to the new extent A1:V6
$table = $ws.ListObjects.Item("contentDataTable")
$table.Resize($ws.Range("A1:V6"))
